$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 282.27777
$ws.Range("I39").Value = 280.1875
$ws.Range("J39").Value = 299
$ws.Range("K39").Value = 840.5625
$ws.Range("L39").Value = 897
$ws.Range("M39").Value = -544.5625
$ws.Range("N39").Value = -1489

# Row 40
$ws.Range("H40").Value = 983.3333
$ws.Range("J40").Value = 1200
$ws.Range("L40").Value = 1200
$ws.Range("N40").Value = -1550

# Row 92
$ws.Range("H92").Value = 2409.182
$ws.Range("I92").Value = 1916.7778
$ws.Range("J92").Value = 4625
$ws.Range("K92").Value = 1916.7778
$ws.Range("L92").Value = 4625
$ws.Range("M92").Value = -668.7778000000001
$ws.Range("N92").Value = -7121

# Row 98
$ws.Range("H98").Value = 4006.1904
$ws.Range("I98").Value = 3924.0715
$ws.Range("K98").Value = 3924.0715
$ws.Range("M98").Value = -2426.0715

# Row 113
$ws.Range("H113").Value = 11829.7
$ws.Range("I113").Value = 17149
$ws.Range("J113").Value = 3850.75
$ws.Range("K113").Value = 17149
$ws.Range("L113").Value = 3850.75
$ws.Range("M113").Value = -13895
$ws.Range("N113").Value = -10358.75

# Row 122
$ws.Range("H122").Value = 4006.1904
$ws.Range("I122").Value = 3924.0715
$ws.Range("K122").Value = 11772.2145
$ws.Range("M122").Value = -9322.2145

# Row 135
$ws.Range("H135").Value = 1851.5
$ws.Range("I135").Value = 1921.8
$ws.Range("K135").Value = 17296.2
$ws.Range("M135").Value = -14761.2

# Row 137
$ws.Range("H137").Value = 16137.8125
$ws.Range("I137").Value = 3149.875
$ws.Range("K137").Value = 9449.625
$ws.Range("M137").Value = -6899.625


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2177332.2
$ws.Range("I32").Value = 2704654.8
$ws.Range("K32").Value = 2704654.8
$ws.Range("M32").Value = -2704367.8

# Row 45
$ws.Range("H45").Value = 1580.5714
$ws.Range("J45").Value = 1560
$ws.Range("L45").Value = 1560
$ws.Range("N45").Value = -2314

# Row 61
$ws.Range("H61").Value = 2869114.8
$ws.Range("I61").Value = 5355
$ws.Range("K61").Value = 5355
$ws.Range("M61").Value = -5143

# Row 110
$ws.Range("H110").Value = 1309
$ws.Range("I110").Value = 1413.5294
$ws.Range("K110").Value = 1413.5294
$ws.Range("M110").Value = 631.4706000000001

# Row 136
$ws.Range("H136").Value = 2869114.8
$ws.Range("I136").Value = 5355
$ws.Range("K136").Value = 16065
$ws.Range("M136").Value = -13515


$ws = $wb.Worksheets.Item("BSM")
# Row 42
$ws.Range("H42").Value = 399684
$ws.Range("J42").Value = 399684
$ws.Range("L42").Value = 399684
$ws.Range("N42").Value = -400340

# Row 86
$ws.Range("H86").Value = 2064
$ws.Range("I86").Value = 2135.0476
$ws.Range("K86").Value = 2135.0476
$ws.Range("M86").Value = -1012.0476

# Row 89
$ws.Range("H89").Value = 2064
$ws.Range("I89").Value = 2135.0476
$ws.Range("K89").Value = 10675.238
$ws.Range("M89").Value = -5059.237999999999


$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 6333.3335
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 6333.3335
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 70
$ws.Range("H70").Value = 17000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73
$ws.Range("H73").Value = 17000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 122
$ws.Range("H122").Value = 2101.8
$ws.Range("I122").Value = 2101.8
$ws.Range("K122").Value = 6305.400000000001
$ws.Range("M122").Value = -3855.400000000001

# Row 132
$ws.Range("H132").Value = 30278840
$ws.Range("I132").Value = 2136.2856
$ws.Range("K132").Value = 6408.8568
$ws.Range("M132").Value = -3878.8568


$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 586.3570999999999
$ws.Range("J68").Value = 714.875
$ws.Range("L68").Value = 2144.625
$ws.Range("N68").Value = -3766.625

# Row 71
$ws.Range("H71").Value = 586.3570999999999
$ws.Range("J71").Value = 714.875
$ws.Range("L71").Value = 6433.875
$ws.Range("N71").Value = -14545.875

# Row 80
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 27500
$ws.Range("L80").Value = 82500
$ws.Range("N80").Value = -84372

# Row 83
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 27500
$ws.Range("L83").Value = 247500
$ws.Range("N83").Value = -256860

# Row 121
$ws.Range("H121").Value = 887.1818
$ws.Range("I121").Value = 477.5
$ws.Range("K121").Value = 1432.5
$ws.Range("M121").Value = -122.5

# Row 129
$ws.Range("H129").Value = 1742.8125
$ws.Range("I129").Value = 1458.8462
$ws.Range("K129").Value = 4376.5386
$ws.Range("M129").Value = 623.4614000000001

# Row 131
$ws.Range("H131").Value = 1485.35
$ws.Range("I131").Value = 1133.3334
$ws.Range("J131").Value = 1496.237
$ws.Range("K131").Value = 3400.0002
$ws.Range("L131").Value = 4488.711
$ws.Range("M131").Value = 1639.9998
$ws.Range("N131").Value = -14568.711


$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 766.75
$ws.Range("I43").Value = 766.75
$ws.Range("K43").Value = 766.75
$ws.Range("M43").Value = -615.75

# Row 113
$ws.Range("H113").Value = 3250.5715
$ws.Range("I113").Value = 2350.8
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 2350.8
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -180.8000000000002
$ws.Range("N113").Value = -9840

# Row 132
$ws.Range("H132").Value = 796072.0600000001
$ws.Range("I132").Value = 5503.4
$ws.Range("K132").Value = 16510.2
$ws.Range("M132").Value = -13980.2


$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4629.6665
$ws.Range("I40").Value = 4443.5
$ws.Range("K40").Value = 4443.5
$ws.Range("M40").Value = -4307.5

# Row 46
$ws.Range("H46").Value = 2746.0908
$ws.Range("I46").Value = 1107.1428
$ws.Range("J46").Value = 3510.9333
$ws.Range("K46").Value = 1107.1428
$ws.Range("L46").Value = 3510.9333
$ws.Range("M46").Value = -919.1428000000001
$ws.Range("N46").Value = -3886.9333

# Row 132
$ws.Range("I132").Value = 2989
$ws.Range("J132").Value = 11637309
$ws.Range("K132").Value = 8967
$ws.Range("L132").Value = 34911927
$ws.Range("M132").Value = -6437
$ws.Range("N132").Value = -34916987

# Row 136
$ws.Range("H136").Value = 165611.92
$ws.Range("I136").Value = 18213.584
$ws.Range("K136").Value = 54640.75199999999
$ws.Range("M136").Value = -52090.75199999999


$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 4931.2144
$ws.Range("I113").Value = 6098.222
$ws.Range("J113").Value = 2830.6
$ws.Range("K113").Value = 18294.666
$ws.Range("L113").Value = 8491.799999999999
$ws.Range("M113").Value = -16124.666
$ws.Range("N113").Value = -12831.8

# Row 135
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140

